$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture a reference style cell (Minion Pro font, used for the 2020 poll rows)
# before we overwrite it, so we can re-apply the same style afterward.
$ws.Range("B11").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- First, overwrite the old 2020-poll label cells (B11:B14) with an already
# existing shared string ("National Exit Poll"). This frees up / recompacts the
# shared-string table slots that used to hold the old long poll names, so that
# the brand new strings we introduce next land at the same indices the
# original authored workbook used. ---
$ws.Range("B11").Value2 = "National Exit Poll"
$ws.Range("B12").Value2 = "National Exit Poll"
$ws.Range("B13").Value2 = "National Exit Poll"
$ws.Range("B14").Value2 = "National Exit Poll"

# These rows no longer belong to the "2020 polls" group, so drop their
# (Minion Pro) font styling back to the plain default style used elsewhere
# in the National Exit Poll rows.
$ws.Range("B5").Copy()
$ws.Range("B11:B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Introduce the brand-new strings in the exact order needed so the shared
# string table ends up ordered the same way as the authored workbook. ---
$ws.Range("B22").Value2 = "The Economist/YouGov"
$ws.Range("B16").Value2 = "WaPo-ABC (March)"
$ws.Range("B17").Value2 = "Pew (Summer)"
$ws.Range("B18").Value2 = "NYT/Siena (June)"
$ws.Range("B19").Value2 = "Marist (Sept.)"
$ws.Range("B20").Value2 = "Quinnipiac (July)"
$ws.Range("B21").Value2 = "Emerson Coll. (July)"
$ws.Range("B23").Value2 = "Nationscape (Apr.-June)"

# --- Rows 11-15: historical National Exit Poll data (moved up from rows 15-19) ---
$ws.Range("A11").Value2 = 2004
$ws.Range("C11").Formula = "=-D11"
$ws.Range("D11").Value2 = 5

$ws.Range("A12").Value2 = 2000
$ws.Range("C12").Formula = "=-D12"
$ws.Range("D12").Value2 = -4

$ws.Range("A13").Value2 = 1996
$ws.Range("C13").Formula = "=-D13"
$ws.Range("D13").Value2 = -6

$ws.Range("A14").Value2 = 1992
$ws.Range("C14").Formula = "=-D14"
$ws.Range("D14").Value2 = -11

$ws.Range("A15").Value2 = 1988
$ws.Range("B15").Value2 = "National Exit Poll"
$ws.Range("C15").Formula = "=-D15"
$ws.Range("D15").Value2 = 2

# --- Rows 16-23: 2020 poll data (shortened labels + two new polls) ---
$ws.Range("A16").Value2 = 2020
$ws.Range("C16").Value2 = 15
$ws.Range("D16").Formula = "=-C16"

$ws.Range("A17").Value2 = 2020
$ws.Range("C17").Value2 = -6
$ws.Range("D17").Formula = "=-C17"

$ws.Range("A18").Value2 = 2020
$ws.Range("C18").Value2 = 2
$ws.Range("D18").Formula = "=-C18"

$ws.Range("A19").Value2 = 2020
$ws.Range("C19").Value2 = 6
$ws.Range("D19").Formula = "=-C19"

$ws.Range("A20").Value2 = 2020
$ws.Range("C20").Value2 = 14
$ws.Range("D20").Formula = "=-C20"

$ws.Range("A21").Value2 = 2020
$ws.Range("C21").Value2 = -20
$ws.Range("D21").Formula = "=-C21"

$ws.Range("A22").Value2 = 2020
$ws.Range("C22").Value2 = 2
$ws.Range("D22").Formula = "=-C22"

$ws.Range("A23").Value2 = 2020
$ws.Range("C23").Value2 = 7.9
$ws.Range("D23").Formula = "=-C23"

# Re-apply the Minion Pro style to the B column for the 2020 poll rows
$ws.Range("Z1").Copy()
$ws.Range("B16:B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# Make column D2:D10 a single fill-down (shared) formula, matching original intent
$ws.Range("D2:D10").Formula = "=-C2"

# Update the selected cell to match the final state
$ws.Range("C23").Select()
